# Add a space between "App" and the number in the five "App<N>" dependency
# names (F2:F6) so they read "App 1" .. "App 5", matching the other
# Dependency_Name groups (e.g. "Person 1", "Tech 1", ...).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "App 1"
$ws.Range("F3").Value = "App 2"
$ws.Range("F4").Value = "App 3"
$ws.Range("F5").Value = "App 4"
$ws.Range("F6").Value = "App 5"

# Move the active selection to F9 (and drop the old C2:C31 selection /
# frozen top-left-cell scroll position) to match the saved view state.
$ws.Range("F9").Select()
